$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(183).Insert()

$ws.Range("A183").Value = 10
$ws.Range("B183").Value = "Vega Modelo de Temuco"
$ws.Range("C183").Value = "La Araucanía"
$ws.Range("D183").Value = 44491
$ws.Range("E183").Value = 9
$ws.Range("F183").Value = 100114013
$ws.Range("G183").Value = "Zanahoria"
$ws.Range("H183").Value = "Sin especificar"
$ws.Range("I183").Value = "Primera"
$ws.Range("J183").Value = 100
$ws.Range("K183").Value = 7000
$ws.Range("L183").Value = 7000
$ws.Range("M183").Value = 7000
$ws.Range("N183").Value = "$/saco 20 kilos"
$ws.Range("O183").Value = "Región del Maule"
$ws.Range("P183").Value = 350
$ws.Range("Q183").Value = 20
$ws.Range("R183").Value = "Hortaliza"

$ws.Range("D183").NumberFormat = "YYYY-MM-DD HH:MM:SS"
